$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the discontinued item row ("ZURCAL 40MG 14 GASTRO RESISTANT TAB") -
# everything below shifts up one row.
$ws.Rows(134).Delete()

# Column A is a plain running serial number (1,2,3,...) independent of the
# product in each row, so after the shift it must stay sequential rather
# than carry up the value from the row below.
for ($r = 134; $r -le 150; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 6
}

# The cached total in column P no longer includes the removed item's price
# (96.00), so update the literal total accordingly.
$ws.Range("P151").Value = 8652.5599999999995

